$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.528.69"
$ws.Cells.Item(2, 5).Value = "  +2.55%  "
$ws.Cells.Item(3, 4).Value = "2.312.08"
$ws.Cells.Item(3, 5).Value = "  +1.58%  "
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
$ws.Cells.Item(5, 4).Value = "'311.14"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.38%  "
$ws.Cells.Item(6, 4).Value = "'102.25"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +3.69%  "
$ws.Cells.Item(7, 5).Value = "  +1.43%  "
$ws.Cells.Item(8, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 4).Value = "'0.531"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +6.96%  "
$ws.Cells.Item(10, 4).Value = "'35.71"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.34%  "
$ws.Cells.Item(11, 5).Value = "  +3.09%  "
$ws.Cells.Item(12, 4).Value = "'0.112"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.51%  "
$ws.Cells.Item(13, 5).Value = "  +1.17%  "
$ws.Cells.Item(14, 4).Value = "2.670.57"
$ws.Cells.Item(14, 5).Value = "  +1.58%  "
$ws.Cells.Item(15, 4).Value = "'15.01"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.40%  "
$ws.Cells.Item(16, 4).Value = "2.316.58"
$ws.Cells.Item(16, 5).Value = "  +1.66%  "
$ws.Cells.Item(17, 4).Value = "'0.809"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.45%  "
$ws.Cells.Item(18, 4).Value = "43.422.08"
$ws.Cells.Item(18, 5).Value = "  +2.69%  "
$ws.Cells.Item(19, 4).Value = "'12.42"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -1.55%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0925"
$ws.Cells.Item(20, 5).Value = "  +1.77%  "
$ws.Cells.Item(21, 4).Value = "'6.18"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +1.97%  "
$ws.Cells.Item(22, 4).Value = "'68.20"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -0.13%  "
$ws.Cells.Item(23, 4).Value = "'241.55"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.99%  "
$ws.Cells.Item(24, 4).Value = "'2.04"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +2.82%  "
$ws.Cells.Item(25, 4).Value = "'2.62"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.71%  "
$ws.Cells.Item(26, 5).Value = "  +0.26%  "
$ws.Cells.Item(27, 5).Value = "  -1.65%  "
$ws.Cells.Item(28, 4).Value = "'24.75"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +4.23%  "
$ws.Cells.Item(29, 2).Value = "Toncoin"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(29, 4).Value = "'2.24"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +5.54%  "
$ws.Cells.Item(30, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(30, 4).Value = "'36.70"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -3.92%  "
$ws.Cells.Item(31, 2).Value = "Cosmos"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(31, 4).Value = "'9.62"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.55%  "
$ws.Cells.Item(32, 4).Value = "'167.89"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +3.65%  "
$ws.Cells.Item(33, 5).Value = "  +0.21%  "
$ws.Cells.Item(34, 5).Value = "  +0.09%  "
$ws.Cells.Item(35, 4).Value = "'2.52"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +6.57%  "
$ws.Cells.Item(37, 4).Value = "'0.0745"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.32%  "
$ws.Cells.Item(38, 5).Value = "  -1.00%  "
$ws.Cells.Item(39, 2).Value = "ARBITRUM"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(39, 4).Value = "'1.87"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +2.11%  "
$ws.Cells.Item(40, 2).Value = "Kaspa"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(40, 4).Value = "'0.106"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +1.10%  "
$ws.Cells.Item(41, 5).Value = "  +1.36%  "
$ws.Cells.Item(42, 4).Value = "'4.28"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +4.35%  "
$ws.Cells.Item(43, 5).Value = "  -0.71%  "
$ws.Cells.Item(44, 4).Value = "'19.31"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +1.99%  "
$ws.Cells.Item(45, 4).Value = "'0.0289"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +1.94%  "
$ws.Cells.Item(46, 4).Value = "1.969.85"
$ws.Cells.Item(46, 5).Value = "  +0.92%  "
$ws.Cells.Item(47, 5).Value = "  +2.05%  "
$ws.Cells.Item(48, 4).Value = "'9.94"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.25%  "
$ws.Cells.Item(49, 4).Value = "'55.51"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +2.92%  "
$ws.Cells.Item(50, 5).Value = "  +6.33%  "
$ws.Cells.Item(51, 4).Value = "'1.57"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +6.31%  "
